$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D column cells keep their text format so numeric-looking strings
# (e.g. "27.319.59", "1.002") are not reinterpreted as numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.319.59"
$ws.Range("E2").Value = "  -0.85%  "
$ws.Range("D3").Value = "1.788.04"
$ws.Range("E3").Value = "  -2.31%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "340.74"
$ws.Range("E5").Value = "  -1.11%  "
$ws.Range("D6").Value = "0.9992"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "0.3961"
$ws.Range("E7").Value = "  +3.30%  "
$ws.Range("D8").Value = "0.3453"
$ws.Range("E8").Value = "  -2.73%  "
$ws.Range("D9").Value = "48.04"
$ws.Range("E9").Value = "  -4.77%  "
$ws.Range("D10").Value = "1.192"
$ws.Range("E10").Value = "  -3.83%  "
$ws.Range("D11").Value = "0.07448"
$ws.Range("E11").Value = "  -4.48%  "
$ws.Range("D12").Value = "0.9995"
$ws.Range("E12").Value = "  -0.15%  "
$ws.Range("D13").Value = "21.68"
$ws.Range("E13").Value = "  -3.16%  "
$ws.Range("D14").Value = "6.452"
$ws.Range("E14").Value = "  -2.62%  "
$ws.Range("D15").Value = "1.787.25"
$ws.Range("E15").Value = "  -2.31%  "
$ws.Range("D16").Value = "7.079"
$ws.Range("E16").Value = "  -2.44%  "
$ws.Range("E17").Value = "  -3.38%  "
$ws.Range("D18").Value = "0.06663"
$ws.Range("E18").Value = "  -1.07%  "
$ws.Range("D19").Value = "83.98"
$ws.Range("E19").Value = "  -3.38%  "
$ws.Range("D20").Value = "0.9989"
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("D21").Value = "17.69"
$ws.Range("E21").Value = "  +0.19%  "
$ws.Range("D22").Value = "6.489"
$ws.Range("E22").Value = "  -1.77%  "
$ws.Range("D23").Value = "27.313.75"
$ws.Range("E23").Value = "  -0.90%  "
$ws.Range("D24").Value = "12.34"
$ws.Range("E24").Value = "  -6.58%  "
$ws.Range("D25").Value = "2.378"
$ws.Range("E25").Value = "  -3.98%  "
$ws.Range("D26").Value = "1.469"
$ws.Range("E26").Value = "  -0.59%  "
$ws.Range("D27").Value = "2.504"
$ws.Range("E27").Value = "  -8.26%  "
$ws.Range("D28").Value = "21.14"
$ws.Range("E28").Value = "  -4.87%  "
$ws.Range("D29").Value = "156.53"
$ws.Range("E29").Value = "  +1.86%  "
$ws.Range("D30").Value = "1.989.92"
$ws.Range("E30").Value = "  -2.25%  "
$ws.Range("D31").Value = "135.36"
$ws.Range("E31").Value = "  -0.36%  "
$ws.Range("D32").Value = "4.014"
$ws.Range("E32").Value = "  -1.53%  "
$ws.Range("D33").Value = "6.008"
$ws.Range("E33").Value = "  -5.96%  "
$ws.Range("D34").Value = "0.08761"
$ws.Range("E34").Value = "  -0.64%  "
$ws.Range("D35").Value = "13.00"
$ws.Range("E35").Value = "  -6.97%  "
$ws.Range("D36").Value = "1.617"
$ws.Range("E36").Value = "  -5.03%  "
$ws.Range("D37").Value = "5.400"
$ws.Range("E37").Value = "  -4.45%  "
$ws.Range("D38").Value = "0.6823"
$ws.Range("E38").Value = "  -3.94%  "
$ws.Range("D39").Value = "0.02380"
$ws.Range("E39").Value = "  -1.59%  "
$ws.Range("D40").Value = "0.06396"
$ws.Range("E40").Value = "  -1.96%  "
$ws.Range("E41").Value = "  -3.01%  "
$ws.Range("D42").Value = "1.246"
$ws.Range("E42").Value = "  -5.06%  "
$ws.Range("D43").Value = "8.418"
$ws.Range("E43").Value = "  -7.89%  "
$ws.Range("D44").Value = "14.38"
$ws.Range("E44").Value = "  -3.02%  "
$ws.Range("D45").Value = "0.9989"
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").Value = "0.6389"
$ws.Range("E46").Value = "  -4.24%  "
$ws.Range("D47").Value = "3.876"
$ws.Range("E47").Value = "  -2.02%  "
$ws.Range("D48").Value = "2.128"
$ws.Range("E48").Value = "  -3.40%  "
$ws.Range("D49").Value = "131.97"
$ws.Range("E49").Value = "  -1.15%  "
$ws.Range("D50").Value = "0.07123"
$ws.Range("E50").Value = "  -2.80%  "
$ws.Range("D51").Value = "78.72"
$ws.Range("E51").Value = "  -3.10%  "
